$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1004.375
$ws.Range("I19").Value = 1110.3636
$ws.Range("K19").Value = 1110.3636
$ws.Range("M19").Value = -935.3635999999999
$ws.Range("H125").Value = 133341840
$ws.Range("I125").Value = 214294420
$ws.Range("J125").Value = 62508336
$ws.Range("K125").Value = 1928649780
$ws.Range("L125").Value = 562575024
$ws.Range("M125").Value = -1928647320
$ws.Range("N125").Value = -562579944

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2791.077
$ws.Range("I61").Value = 2791.077
$ws.Range("K61").Value = 2791.077
$ws.Range("M61").Value = -2579.077
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H122").Value = 2028.375
$ws.Range("I122").Value = 806.44446
$ws.Range("K122").Value = 2419.33338
$ws.Range("M122").Value = 30.66661999999997
$ws.Range("H132").Value = 1510.4
$ws.Range("I132").Value = 1510.4
$ws.Range("K132").Value = 4531.200000000001
$ws.Range("M132").Value = -2001.200000000001
$ws.Range("H136").Value = 2791.077
$ws.Range("I136").Value = 2791.077
$ws.Range("K136").Value = 8373.231
$ws.Range("M136").Value = -5823.231

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 35282.285
$ws.Range("J81").Value = 38629.332
$ws.Range("L81").Value = 38629.332
$ws.Range("N81").Value = -40751.332
$ws.Range("H84").Value = 35282.285
$ws.Range("J84").Value = 38629.332
$ws.Range("L84").Value = 115887.996
$ws.Range("N84").Value = -126495.996
$ws.Range("H134").Value = 2088.6
$ws.Range("I134").Value = 2088.6
$ws.Range("K134").Value = 6265.799999999999
$ws.Range("M134").Value = -3730.799999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2383.5
$ws.Range("I31").Value = 1796.6
$ws.Range("J31").Value = 3361.6667
$ws.Range("K31").Value = 1796.6
$ws.Range("L31").Value = 3361.6667
$ws.Range("M31").Value = -1501.6
$ws.Range("N31").Value = -3951.6667
$ws.Range("H34").Value = 2383.5
$ws.Range("I34").Value = 1796.6
$ws.Range("J34").Value = 3361.6667
$ws.Range("K34").Value = 1796.6
$ws.Range("L34").Value = 3361.6667
$ws.Range("M34").Value = -1594.6
$ws.Range("N34").Value = -3765.6667
$ws.Range("H58").Value = 1712.4348
$ws.Range("I58").Value = 1800.7894
$ws.Range("K58").Value = 1800.7894
$ws.Range("M58").Value = -1597.7894
$ws.Range("H93").Value = 6333
$ws.Range("I93").Value = 6333
$ws.Range("K93").Value = 6333
$ws.Range("M93").Value = -4461
$ws.Range("H99").Value = 3602.8572
$ws.Range("I99").Value = 2110.5
$ws.Range("J99").Value = 4199.8
$ws.Range("K99").Value = 2110.5
$ws.Range("L99").Value = 4199.8
$ws.Range("M99").Value = -612.5
$ws.Range("N99").Value = -7195.8
$ws.Range("H122").Value = 1476.8
$ws.Range("I122").Value = 1501.1578
$ws.Range("K122").Value = 4503.4734
$ws.Range("M122").Value = -2053.4734
$ws.Range("H126").Value = 3602.8572
$ws.Range("I126").Value = 2110.5
$ws.Range("J126").Value = 4199.8
$ws.Range("K126").Value = 6331.5
$ws.Range("L126").Value = 12599.4
$ws.Range("M126").Value = -3861.5
$ws.Range("N126").Value = -17539.4
$ws.Range("H132").Value = 2909.7273
$ws.Range("I132").Value = 3063.375
$ws.Range("K132").Value = 9190.125
$ws.Range("M132").Value = -6660.125
$ws.Range("H134").Value = 904.7
$ws.Range("I134").Value = 904.7
$ws.Range("K134").Value = 2714.1
$ws.Range("M134").Value = -179.1000000000004
$ws.Range("H136").Value = 1712.4348
$ws.Range("I136").Value = 1800.7894
$ws.Range("K136").Value = 5402.3682
$ws.Range("M136").Value = -2852.3682

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H37").Value = 59968.332
$ws.Range("J37").Value = 59968.332
$ws.Range("L37").Value = 179904.996
$ws.Range("N37").Value = -180128.996
$ws.Range("H64").Value = 8999
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 8999
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 26997
$ws.Range("M64").ClearContents()
$ws.Range("N64").Value = -27537
$ws.Range("H67").Value = 8999
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 8999
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 26997
$ws.Range("M67").ClearContents()
$ws.Range("N67").Value = -28869
$ws.Range("H98").Value = 2086.5715
$ws.Range("I98").Value = 650
$ws.Range("J98").Value = 4002
$ws.Range("K98").Value = 1950
$ws.Range("L98").Value = 12006
$ws.Range("M98").Value = -452
$ws.Range("N98").Value = -15002
$ws.Range("H117").Value = 10361.818
$ws.Range("I117").Value = 391.6
$ws.Range("J117").Value = 18670.334
$ws.Range("K117").Value = 1174.8
$ws.Range("L117").Value = 56011.00199999999
$ws.Range("M117").Value = 2267.2
$ws.Range("N117").Value = -62895.00199999999
$ws.Range("H129").Value = 1252723.9
$ws.Range("I129").Value = 525
$ws.Range("J129").Value = 2004043.2
$ws.Range("K129").Value = 1575
$ws.Range("L129").Value = 6012129.6
$ws.Range("M129").Value = 3425
$ws.Range("N129").Value = -6022129.6
$ws.Range("H131").Value = 305390
$ws.Range("J131").Value = 419381
$ws.Range("L131").Value = 1258143
$ws.Range("N131").Value = -1268223
$ws.Range("H137").Value = 2949.6667
$ws.Range("I137").Value = 1824.75
$ws.Range("J137").Value = 5199.5
$ws.Range("K137").Value = 5474.25
$ws.Range("L137").Value = 15598.5
$ws.Range("M137").Value = -374.25
$ws.Range("N137").Value = -25798.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 833
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 499
$ws.Range("K113").Value = 1000
$ws.Range("L113").Value = 499
$ws.Range("M113").Value = 1170
$ws.Range("N113").Value = -4839

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7693.8887
$ws.Range("I7").Value = 4500
$ws.Range("J7").Value = 8922.308000000001
$ws.Range("K7").Value = 4500
$ws.Range("L7").Value = 8922.308000000001
$ws.Range("M7").Value = -4388
$ws.Range("N7").Value = -9146.308000000001
$ws.Range("H55").Value = 198
$ws.Range("J55").Value = 299.5
$ws.Range("L55").Value = 299.5
$ws.Range("N55").Value = -645.5
$ws.Range("H61").Value = 3278.3333
$ws.Range("I61").Value = 3425
$ws.Range("J61").Value = 2985
$ws.Range("K61").Value = 3425
$ws.Range("L61").Value = 2985
$ws.Range("M61").Value = -3223
$ws.Range("N61").Value = -3389
$ws.Range("H113").Value = 3278.3333
$ws.Range("I113").Value = 3425
$ws.Range("J113").Value = 2985
$ws.Range("K113").Value = 3425
$ws.Range("L113").Value = 2985
$ws.Range("M113").Value = -1255
$ws.Range("N113").Value = -7325
$ws.Range("H122").Value = 6112.375
$ws.Range("I122").Value = 7634.3076
$ws.Range("J122").Value = 5071.0527
$ws.Range("K122").Value = 22902.9228
$ws.Range("L122").Value = 15213.1581
$ws.Range("M122").Value = -20452.9228
$ws.Range("N122").Value = -20113.1581
$ws.Range("H126").Value = 7693.8887
$ws.Range("I126").Value = 4500
$ws.Range("J126").Value = 8922.308000000001
$ws.Range("K126").Value = 13500
$ws.Range("L126").Value = 26766.924
$ws.Range("M126").Value = -11030
$ws.Range("N126").Value = -31706.924
$ws.Range("H136").Value = 3757.5
$ws.Range("I136").Value = 2911.6667
$ws.Range("J136").Value = 6295
$ws.Range("K136").Value = 8735.000100000001
$ws.Range("L136").Value = 18885
$ws.Range("M136").Value = -6185.000100000001
$ws.Range("N136").Value = -23985

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 883.9
$ws.Range("I122").Value = 755
$ws.Range("K122").Value = 2265
$ws.Range("M122").Value = 185
$ws.Range("H136").Value = 3124.25
$ws.Range("I136").Value = 3104.4736
$ws.Range("K136").Value = 9313.4208
$ws.Range("M136").Value = -6763.4208
